# Auto-generated script to apply market-data refresh to Kujata Profits workbook
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per-row
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 113.95652
$ws.Range("I33").Value = 92.47619
$ws.Range("K33").Value = 92.47619
$ws.Range("M33").Value = 136.52381
$ws.Range("H86").Value = 4213.9
$ws.Range("I86").Value = 4247.25
$ws.Range("J86").Value = 4191.6665
$ws.Range("K86").Value = 4247.25
$ws.Range("L86").Value = 4191.6665
$ws.Range("M86").Value = -3124.25
$ws.Range("N86").Value = -6437.6665
$ws.Range("H87").Value = 42749.75
$ws.Range("J87").Value = 42749.75
$ws.Range("L87").Value = 42749.75
$ws.Range("N87").Value = -45245.75
$ws.Range("H88").Value = 2469716.5
$ws.Range("J88").Value = 6172989.5
$ws.Range("L88").Value = 6172989.5
$ws.Range("N88").Value = -6173801.5
$ws.Range("H89").Value = 4213.9
$ws.Range("I89").Value = 4247.25
$ws.Range("J89").Value = 4191.6665
$ws.Range("K89").Value = 21236.25
$ws.Range("L89").Value = 20958.3325
$ws.Range("M89").Value = -15620.25
$ws.Range("N89").Value = -32190.3325
$ws.Range("H90").Value = 42749.75
$ws.Range("J90").Value = 42749.75
$ws.Range("L90").Value = 128249.25
$ws.Range("N90").Value = -140729.25
$ws.Range("H91").Value = 2469716.5
$ws.Range("J91").Value = 6172989.5
$ws.Range("L91").Value = 6172989.5
$ws.Range("N91").Value = -6175797.5
$ws.Range("H92").Value = 759.1
$ws.Range("I92").Value = 759.1
$ws.Range("K92").Value = 759.1
$ws.Range("M92").Value = 488.9
$ws.Range("H116").Value = 3540.4
$ws.Range("I116").Value = 2827.2727
$ws.Range("K116").Value = 2827.2727
$ws.Range("M116").Value = 614.7273
$ws.Range("H137").Value = 1225.6129
$ws.Range("I137").Value = 1090.5264
$ws.Range("K137").Value = 3271.5792
$ws.Range("M137").Value = -721.5792000000001
$ws.Range("H138").Value = 2106.0505
$ws.Range("J138").Value = 2303.9058
$ws.Range("L138").Value = 6911.7174
$ws.Range("N138").Value = -17191.7174

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2011.2
$ws.Range("I45").Value = 1730.2858
$ws.Range("J45").Value = 2666.6667
$ws.Range("K45").Value = 1730.2858
$ws.Range("L45").Value = 2666.6667
$ws.Range("M45").Value = -1353.2858
$ws.Range("N45").Value = -3420.6667
$ws.Range("H61").Value = 1377.5333
$ws.Range("I61").Value = 1212.5385
$ws.Range("J61").Value = 2450
$ws.Range("K61").Value = 1212.5385
$ws.Range("L61").Value = 2450
$ws.Range("M61").Value = -1000.5385
$ws.Range("N61").Value = -2874
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("N107").Value = 0
$ws.Range("H119").Value = 35698
$ws.Range("J119").Value = 35698
$ws.Range("L119").Value = 35698
$ws.Range("N119").Value = -45374
$ws.Range("H132").Value = 2741.3225
$ws.Range("I132").Value = 2519.76
$ws.Range("K132").Value = 7559.280000000001
$ws.Range("M132").Value = -5029.280000000001
$ws.Range("H136").Value = 1377.5333
$ws.Range("I136").Value = 1212.5385
$ws.Range("J136").Value = 2450
$ws.Range("K136").Value = 3637.6155
$ws.Range("L136").Value = 7350
$ws.Range("M136").Value = -1087.6155
$ws.Range("N136").Value = -12450

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3895.389
$ws.Range("I86").Value = 3724.5293
$ws.Range("K86").Value = 3724.5293
$ws.Range("M86").Value = -2601.5293
$ws.Range("H89").Value = 3895.389
$ws.Range("I89").Value = 3724.5293
$ws.Range("K89").Value = 18622.6465
$ws.Range("M89").Value = -13006.6465
$ws.Range("H94").Value = 250000000
$ws.Range("I94").Value = 250000000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 250000000
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -249999549
$ws.Range("H107").Value = 1973.909
$ws.Range("I107").Value = 1480
$ws.Range("J107").Value = 2385.5
$ws.Range("K107").Value = 1480
$ws.Range("L107").Value = 2385.5
$ws.Range("M107").Value = 440
$ws.Range("N107").Value = -6225.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 369.6
$ws.Range("I22").Value = 337.5
$ws.Range("J22").Value = 498
$ws.Range("K22").Value = 337.5
$ws.Range("L22").Value = 498
$ws.Range("M22").Value = 12.5
$ws.Range("N22").Value = -1198
$ws.Range("H31").Value = 815.4154
$ws.Range("I31").Value = 746.7857
$ws.Range("J31").Value = 940.73914
$ws.Range("K31").Value = 746.7857
$ws.Range("L31").Value = 940.73914
$ws.Range("M31").Value = -451.7857
$ws.Range("N31").Value = -1530.73914
$ws.Range("H34").Value = 815.4154
$ws.Range("I34").Value = 746.7857
$ws.Range("J34").Value = 940.73914
$ws.Range("K34").Value = 746.7857
$ws.Range("L34").Value = 940.73914
$ws.Range("M34").Value = -544.7857
$ws.Range("N34").Value = -1344.73914
$ws.Range("H86").Value = 2909181.8
$ws.Range("I86").Value = 5557488
$ws.Range("J86").Value = 20120.637
$ws.Range("K86").Value = 5557488
$ws.Range("L86").Value = 20120.637
$ws.Range("M86").Value = -5556365
$ws.Range("N86").Value = -22366.637
$ws.Range("H89").Value = 2909181.8
$ws.Range("I89").Value = 5557488
$ws.Range("J89").Value = 20120.637
$ws.Range("K89").Value = 27787440
$ws.Range("L89").Value = 100603.185
$ws.Range("M89").Value = -27781824
$ws.Range("N89").Value = -111835.185
$ws.Range("H99").Value = 2002.9166
$ws.Range("I99").Value = 1979.5
$ws.Range("J99").Value = 2049.75
$ws.Range("K99").Value = 1979.5
$ws.Range("L99").Value = 2049.75
$ws.Range("M99").Value = -481.5
$ws.Range("N99").Value = -5045.75
$ws.Range("H105").Value = 1100
$ws.Range("I105").Value = 1000
$ws.Range("K105").Value = 1000
$ws.Range("M105").Value = 747
$ws.Range("H126").Value = 2002.9166
$ws.Range("I126").Value = 1979.5
$ws.Range("J126").Value = 2049.75
$ws.Range("K126").Value = 5938.5
$ws.Range("L126").Value = 6149.25
$ws.Range("M126").Value = -3468.5
$ws.Range("N126").Value = -11089.25
$ws.Range("H132").Value = 6195.727
$ws.Range("I132").Value = 7123.9414
$ws.Range("J132").Value = 3039.8
$ws.Range("K132").Value = 21371.8242
$ws.Range("L132").Value = 9119.400000000001
$ws.Range("M132").Value = -18841.8242
$ws.Range("N132").Value = -14179.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 125162.375
$ws.Range("I8").Value = 125162.375
$ws.Range("K8").Value = 375487.125
$ws.Range("M8").Value = -375348.125
$ws.Range("H20").Value = 1000
$ws.Range("J20").Value = 1000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3454
$ws.Range("H107").Value = 4899.913
$ws.Range("J107").Value = 11705.444
$ws.Range("L107").Value = 35116.33199999999
$ws.Range("N107").Value = -38956.33199999999
$ws.Range("H131").Value = 40001820
$ws.Range("J131").Value = 2249.6667
$ws.Range("L131").Value = 6749.000100000001
$ws.Range("N131").Value = -16829.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226
$ws.Range("H126").Value = 2503.7144
$ws.Range("I126").Value = 1781.5
$ws.Range("K126").Value = 5344.5
$ws.Range("M126").Value = -2874.5
$ws.Range("H132").Value = 2006.5151
$ws.Range("I132").Value = 1622.1111
$ws.Range("J132").Value = 3736.3333
$ws.Range("K132").Value = 4866.3333
$ws.Range("L132").Value = 11208.9999
$ws.Range("M132").Value = -2336.3333
$ws.Range("N132").Value = -16268.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1399.8572
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 1999.75
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 1999.75
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -2589.75
$ws.Range("H27").Value = 1399.8572
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 1999.75
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 1999.75
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -2213.75
$ws.Range("H40").Value = 3027
$ws.Range("I40").Value = 2547.5
$ws.Range("J40").Value = 3666.3333
$ws.Range("K40").Value = 2547.5
$ws.Range("L40").Value = 3666.3333
$ws.Range("M40").Value = -2411.5
$ws.Range("N40").Value = -3938.3333
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = 0
$ws.Range("H132").Value = 50133.285
$ws.Range("I132").Value = 2355.3333
$ws.Range("J132").Value = 85966.75
$ws.Range("K132").Value = 7065.999899999999
$ws.Range("L132").Value = 257900.25
$ws.Range("M132").Value = -4535.999899999999
$ws.Range("N132").Value = -262960.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("N57").Value = 0
$ws.Range("H96").Value = 1354.1
$ws.Range("I96").Value = 1282.4445
$ws.Range("K96").Value = 1282.4445
$ws.Range("M96").Value = 90.55549999999994
$ws.Range("H126").Value = 111112790
$ws.Range("I126").Value = 200000930
$ws.Range("J126").Value = 2635
$ws.Range("K126").Value = 600002790
$ws.Range("L126").Value = 7905
$ws.Range("M126").Value = -600000320
$ws.Range("N126").Value = -12845

